$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.32"
$ws.Range("E2").Value = "'-2.11%"
$ws.Range("D3").Value = "'37.78"
$ws.Range("E3").Value = "'-4.32%"
$ws.Range("D4").Value = "'5.040"
$ws.Range("E4").Value = "'-2.15%"
$ws.Range("D5").Value = "'0.07897"
$ws.Range("E5").Value = "'-3.26%"
$ws.Range("D6").Value = "'2.051"
$ws.Range("E6").Value = "'3.73%"
$ws.Range("D7").Value = "'4.391"
$ws.Range("E7").Value = "'3.46%"
$ws.Range("D8").Value = "'8.232"
$ws.Range("E8").Value = "'0.40%"
$ws.Range("D9").Value = "'3.084"
$ws.Range("E9").Value = "'-4.60%"
$ws.Range("D10").Value = "'0.9272"
$ws.Range("E10").Value = "'0.24%"
$ws.Range("D11").Value = "'0.1288"
$ws.Range("E11").Value = "'-8.18%"
$ws.Range("E12").Value = "'-3.76%"
$ws.Range("D13").Value = "'0.08702"
$ws.Range("E13").Value = "'-3.58%"
$ws.Range("D14").Value = "'0.03457"
$ws.Range("E14").Value = "'-0.97%"
$ws.Range("D15").Value = "'0.09746"
$ws.Range("E15").Value = "'-0.42%"
$ws.Range("D16").Value = "'0.001391"
$ws.Range("E16").Value = "'-2.65%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006149"
$ws.Range("E17").Value = "'1.75%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.549"
$ws.Range("E18").Value = "'-2.80%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3442"
$ws.Range("E19").Value = "'-0.38%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1301"
$ws.Range("E20").Value = "'-3.33%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'5.015"
$ws.Range("E21").Value = "'7.16%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2517"
$ws.Range("E22").Value = "'3.88%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04329"
$ws.Range("E23").Value = "'-1.08%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'-0.34%"
$ws.Range("D25").Value = "'0.004601"
$ws.Range("E25").Value = "'-3.90%"
$ws.Range("E26").Value = "'177.14%"
$ws.Range("D39").Value = "'0.02247"
$ws.Range("E39").Value = "'3.62%"
$ws.Range("D40").Value = "'0.05066"
$ws.Range("E40").Value = "'-2.75%"
$ws.Range("D41").Value = "'0.007500"
$ws.Range("E41").Value = "'-0.69%"
$ws.Range("D42").Value = "'0.009938"
$ws.Range("E42").Value = "'0.99%"
$ws.Range("D43").Value = "'0.1359"
$ws.Range("E43").Value = "'-0.87%"
$ws.Range("E44").Value = "'-0.68%"
$ws.Range("D45").Value = "'0.008818"
$ws.Range("E45").Value = "'-1.83%"
$ws.Range("D46").Value = "'0.00006497"
$ws.Range("E46").Value = "'1.61%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.49%"
$ws.Range("D48").Value = "'0.003005"
$ws.Range("E48").Value = "'8.87%"
$ws.Range("D49").Value = "'0.001203"
$ws.Range("E49").Value = "'0.49%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.49%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.49%"
